# Update the "取得日時" (retrieved-at) timestamp in column A for rows 2-12
# of the "ランサーズ" sheet from 2025-10-08 06:27:33 to 2025-10-08 06:33:58
# (commit message: "Append: 2025-10-08 06:34 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-08 06:33:58"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
